# Aspect bar for search SO Library
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the existing single-cell formulas in row 23 (C23:N23) and
# row 26 (C26:N26) into shared formulas by re-entering the formula
# across the whole range at once (mirrors Excel's fill-right/shared
# formula behaviour).
$ws.Range("C23:N23").Formula = "=C16/C15"
$ws.Range("C26:N26").Formula = "=C16/C11"

# New "aspect bar" (T / C / Score) helper table below the report.
$ws.Range("B32").Value = "minimize"
$ws.Range("B33").Value = "maximize"

$ws.Range("C31").Value = "T"
$ws.Range("D31").Value = "C"
$ws.Range("E31").Value = "Score"

$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = 20
$ws.Range("E32").Formula = "=(C32/D32)*100"

$ws.Range("C33").Value = 10
$ws.Range("D33").Value = 11
$ws.Range("E33").Formula = "=(D33/C33)*100"

$ws.Range("D36").Value = "C tidak boleh 0"

# Update the view so it reflects the newly added rows being in focus.
$ws.Range("G30").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
